$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "(308051846, Eyal  Sofer: -9,-5)"
$ws.Range("B1").Value = "(312049950, Molham  Peretz: -8,-7)"
$ws.Range("C1").Value = "(308073899, Anan  Kirshenbaum: -10,-7)"
$ws.Range("D1").Value = "(318869187, Soaad  Leibovich: -10,0)"
$ws.Range("E1").Value = "(205898513, Asaf  Braymok: 7,-1)"
$ws.Range("F1").Value = "(318428158, Tal  Asulin: -4,4)"
$ws.Range("G1").Value = "(316028364, Sami  Castro: 8,0)"

$ws.Range("A3").Value = "cost: 358.6424165715349"
$ws.Range("A4").Value = "time: 66.72848331430698"
